# Revert "Powerpoint writer: consolidate text run nodes."
#
# Splits the trailing-space-joined title runs back into separate
# "word" / " " runs on the three title placeholders that were
# previously consolidated (slide 1: "Header " + "with ", slide 2:
# "Syntax ", slide 3: "Two " + "column "). Re-assigning `.Text` on a
# `Characters(start, length)` sub-range rewrites just that span in
# place, splitting the backing run without touching neighboring runs
# or their formatting.

$p = $ppt.ActivePresentation

# Slide 1 title: "Header with inline code" -> "Header" / " " / "with" / " " / "inline code"
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 6).Text = "Header"
$tr1.Characters(7, 1).Text = " "
$tr1.Characters(8, 4).Text = "with"
$tr1.Characters(12, 1).Text = " "

# Slide 2 title: "Syntax highlighting" -> "Syntax" / " " / "highlighting"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, 6).Text = "Syntax"
$tr2.Characters(7, 1).Text = " "

# Slide 3 title: "Two column slide" -> "Two" / " " / "column" / " " / "slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, 3).Text = "Two"
$tr3.Characters(4, 1).Text = " "
$tr3.Characters(5, 6).Text = "column"
$tr3.Characters(11, 1).Text = " "
